$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value2 = 0.295620437956204
$ws.Range("C2").Value2 = 0.29683698296837
$ws.Range("D2").Value2 = 0.152068126520681
$ws.Range("E2").Value2 = 0.0693430656934307
$ws.Range("F2").Value2 = 0.346715328467153
$ws.Range("G2").Value2 = 0.217761557177616
$ws.Range("H2").Value2 = 0.237226277372263
$ws.Range("I2").Value2 = 0.13625304136253
$ws.Range("J2").Value2 = 0.0559610705596107
$ws.Range("K2").Value2 = 0.0279805352798054
$ws.Range("L2").Value2 = 0.0705596107055961
$ws.Range("M2").Value2 = 0.183698296836983
$ws.Range("N2").Value2 = 0.0340632603406326
$ws.Range("O2").Value2 = 0.192214111922141
$ws.Range("P2").Value2 = 0.531630170316302
$ws.Range("Q2").Value2 = 0.0656934306569343
$ws.Range("R2").Value2 = 0.01338199513382
$ws.Range("S2").Value2 = 0.21654501216545
$ws.Range("T2").Value2 = 0.327250608272506
$ws.Range("U2").Value2 = 0.265206812652068
$ws.Range("W2").Value2 = 0.416058394160584
$ws.Range("X2").Value2 = 0.115571776155718

# Row 3
$ws.Range("B3").Value2 = 0.262773722627737
$ws.Range("C3").Value2 = 0.138686131386861
$ws.Range("D3").Value2 = 0.0316301703163017
$ws.Range("E3").Value2 = 0.0182481751824818
$ws.Range("F3").Value2 = 0.165450121654501
$ws.Range("G3").Value2 = 0.0462287104622871
$ws.Range("H3").Value2 = 0.36374695863747
$ws.Range("I3").Value2 = 0.0583941605839416
$ws.Range("J3").Value2 = 0.759124087591241
$ws.Range("K3").Value2 = 0.347931873479319
$ws.Range("L3").Value2 = 0.571776155717762
$ws.Range("M3").Value2 = 0.0425790754257908
$ws.Range("N3").Value2 = 0.375912408759124
$ws.Range("O3").Value2 = 0.355231143552311
$ws.Range("P3").Value2 = 0.0608272506082725
$ws.Range("Q3").Value2 = 0.58029197080292
$ws.Range("R3").Value2 = 0.0425790754257908
$ws.Range("S3").Value2 = 0.395377128953771
$ws.Range("T3").Value2 = 0.0462287104622871
$ws.Range("U3").Value2 = 0.0072992700729927
$ws.Range("W3").Value2 = 0.193430656934307
$ws.Range("X3").Value2 = 0.386861313868613

# Row 4
$ws.Range("B4").Value2 = 0.413625304136253
$ws.Range("C4").Value2 = 0.312652068126521
$ws.Range("D4").Value2 = 0.498783454987835
$ws.Range("E4").Value2 = 0.5669099756691
$ws.Range("F4").Value2 = 0.232360097323601
$ws.Range("G4").Value2 = 0.591240875912409
$ws.Range("H4").Value2 = 0.114355231143552
$ws.Range("I4").Value2 = 0.250608272506083
$ws.Range("J4").Value2 = 0.0109489051094891
$ws.Range("K4").Value2 = 0.171532846715328
$ws.Range("L4").Value2 = 0.190997566909976
$ws.Range("M4").Value2 = 0.512165450121654
$ws.Range("N4").Value2 = 0.0279805352798054
$ws.Range("O4").Value2 = 0.354014598540146
$ws.Range("P4").Value2 = 0.0875912408759124
$ws.Range("Q4").Value2 = 0.263990267639903
$ws.Range("R4").Value2 = 0.341849148418492
$ws.Range("S4").Value2 = 0.0267639902676399
$ws.Range("T4").Value2 = 0.403892944038929
$ws.Range("U4").Value2 = 0.367396593673966
$ws.Range("W4").Value2 = 0.0072992700729927
$ws.Range("X4").Value2 = 0.25669099756691

# Row 5
$ws.Range("B5").Value2 = 0.0279805352798054
$ws.Range("C5").Value2 = 0.251824817518248
$ws.Range("D5").Value2 = 0.317518248175182
$ws.Range("E5").Value2 = 0.345498783454988
$ws.Range("F5").Value2 = 0.255474452554745
$ws.Range("G5").Value2 = 0.144768856447689
$ws.Range("H5").Value2 = 0.284671532846715
$ws.Range("I5").Value2 = 0.554744525547445
$ws.Range("J5").Value2 = 0.173965936739659
$ws.Range("K5").Value2 = 0.452554744525547
$ws.Range("L5").Value2 = 0.166666666666667
$ws.Range("M5").Value2 = 0.261557177615572
$ws.Range("N5").Value2 = 0.562043795620438
$ws.Range("O5").Value2 = 0.0985401459854015
$ws.Range("P5").Value2 = 0.319951338199513
$ws.Range("Q5").Value2 = 0.0900243309002433
$ws.Range("R5").Value2 = 0.602189781021898
$ws.Range("S5").Value2 = 0.361313868613139
$ws.Range("T5").Value2 = 0.222627737226277
$ws.Range("U5").Value2 = 0.360097323600973
$ws.Range("W5").Value2 = 0.383211678832117
$ws.Range("X5").Value2 = 0.240875912408759
